# Update the "maximo" / time metric and re-shuffle the Solucion (Pedido/Salida)
# assignment to reflect a new randomized heuristic run.

$wb = $excel.ActiveWorkbook

$wsResumen  = $wb.Worksheets.Item("Resumen")
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

# --- Resumen sheet: updated Maximo value for zone Z1 ---
$wsResumen.Range("C2").Value = 644.4108011676939

# --- Solucion sheet: new Pedido -> Salida assignment (rows 2-41) ---
$solucion = @(
    @("Pedido_23", "S001"),
    @("Pedido_15", "S025"),
    @("Pedido_20", "S005"),
    @("Pedido_27", "S029"),
    @("Pedido_30", "S002"),
    @("Pedido_35", "S006"),
    @("Pedido_5",  "S026"),
    @("Pedido_36", "S030"),
    @("Pedido_7",  "S003"),
    @("Pedido_19", "S027"),
    @("Pedido_12", "S007"),
    @("Pedido_22", "S031"),
    @("Pedido_4",  "S004"),
    @("Pedido_25", "S028"),
    @("Pedido_18", "S008"),
    @("Pedido_13", "S032"),
    @("Pedido_9",  "S009"),
    @("Pedido_11", "S033"),
    @("Pedido_3",  "S013"),
    @("Pedido_40", "S037"),
    @("Pedido_38", "S010"),
    @("Pedido_33", "S034"),
    @("Pedido_17", "S014"),
    @("Pedido_6",  "S038"),
    @("Pedido_16", "S011"),
    @("Pedido_10", "S035"),
    @("Pedido_39", "S015"),
    @("Pedido_1",  "S039"),
    @("Pedido_31", "S012"),
    @("Pedido_14", "S036"),
    @("Pedido_28", "S016"),
    @("Pedido_2",  "S040"),
    @("Pedido_32", "S017"),
    @("Pedido_37", "S021"),
    @("Pedido_34", "S018"),
    @("Pedido_24", "S022"),
    @("Pedido_21", "S019"),
    @("Pedido_8",  "S023"),
    @("Pedido_29", "S020"),
    @("Pedido_26", "S024")
)

for ($i = 0; $i -lt $solucion.Length; $i++) {
    $row = $i + 2
    $wsSolucion.Cells.Item($row, 1).Value = $solucion[$i][0]
    $wsSolucion.Cells.Item($row, 2).Value = $solucion[$i][1]
}

# --- Metricas sheet: updated Tiempo values for Z1 and Z2 ---
$wsMetricas.Range("B2").Value = 644.4108011676939
$wsMetricas.Range("B3").Value = 469.6989674559413
